# Update the NATMI LR-pair table (Vtn-Itga8) with refreshed TPM-derived
# statistics: a new "Resolving-Mac" cluster is added (rows 17-21 become a
# full new block) and all existing rows 2-16 get recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Itga8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.337313
$ws.Range("H2").Value = 10.011939
$ws.Range("I2").Value = 0.1958858017947999
$ws.Range("J2").Value = 0.1958858017947999
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3708513333333334
$ws.Range("N2").Value = 1.112554
$ws.Range("O2").Value = 0.07446453063452614
$ws.Range("P2").Value = 0.07446453063452615
$ws.Range("Q2").Value = 1.237646975800667
$ws.Range("R2").Value = 11.138822782206
$ws.Range("S2").Value = 0.01458654428861759
$ws.Range("T2").Value = 0.0145865442886176

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Itga8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.337313
$ws.Range("H3").Value = 10.011939
$ws.Range("I3").Value = 0.1958858017947999
$ws.Range("J3").Value = 0.1958858017947999
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.230730666666667
$ws.Range("N3").Value = 6.692192
$ws.Range("O3").Value = 0.4479161786269527
$ws.Range("P3").Value = 0.4479161786269528
$ws.Range("Q3").Value = 7.444646453365333
$ws.Range("R3").Value = 67.001818080288
$ws.Range("S3").Value = 0.08774041978720346
$ws.Range("T3").Value = 0.08774041978720347

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Itga8"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.337313
$ws.Range("H4").Value = 10.011939
$ws.Range("I4").Value = 0.1958858017947999
$ws.Range("J4").Value = 0.1958858017947999
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.029162
$ws.Range("N4").Value = 3.087486
$ws.Range("O4").Value = 0.2066490218278579
$ws.Range("P4").Value = 0.2066490218278579
$ws.Range("Q4").Value = 3.434635721706
$ws.Range("R4").Value = 30.911721495354
$ws.Range("S4").Value = 0.04047960933086105
$ws.Range("T4").Value = 0.04047960933086105

# Row 5: ECs -> MuSCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Itga8"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.337313
$ws.Range("H5").Value = 10.011939
$ws.Range("I5").Value = 0.1958858017947999
$ws.Range("J5").Value = 0.1958858017947999
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8350726666666667
$ws.Range("N5").Value = 2.505218
$ws.Range("O5").Value = 0.1676771487111334
$ws.Range("P5").Value = 0.1676771487111334
$ws.Range("Q5").Value = 2.786898866411333
$ws.Range("R5").Value = 25.082089797702
$ws.Range("S5").Value = 0.03284557271794627
$ws.Range("T5").Value = 0.03284557271794628

# Row 6: ECs -> Resolving-Mac
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Itga8"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.337313
$ws.Range("H6").Value = 10.011939
$ws.Range("I6").Value = 0.1958858017947999
$ws.Range("J6").Value = 0.1958858017947999
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5144246666666666
$ws.Range("N6").Value = 1.543274
$ws.Range("O6").Value = 0.1032931201995298
$ws.Range("P6").Value = 0.1032931201995298
$ws.Range("Q6").Value = 1.716796127587333
$ws.Range("R6").Value = 15.451165148286
$ws.Range("S6").Value = 0.02023365567017154
$ws.Range("T6").Value = 0.02023365567017154

# Row 7: FAPs -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Itga8"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.28369633333333
$ws.Range("H7").Value = 30.851089
$ws.Range("I7").Value = 0.6036083824529627
$ws.Range("J7").Value = 0.6036083824529627
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3708513333333334
$ws.Range("N7").Value = 1.112554
$ws.Range("O7").Value = 0.07446453063452614
$ws.Range("P7").Value = 0.07446453063452615
$ws.Range("Q7").Value = 3.813722496811778
$ws.Range("R7").Value = 34.323502471306
$ws.Range("S7").Value = 0.04494741488642541
$ws.Range("T7").Value = 0.04494741488642542

# Row 8: FAPs -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Itga8"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.28369633333333
$ws.Range("H8").Value = 30.851089
$ws.Range("I8").Value = 0.6036083824529627
$ws.Range("J8").Value = 0.6036083824529627
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.230730666666667
$ws.Range("N8").Value = 6.692192
$ws.Range("O8").Value = 0.4479161786269527
$ws.Range("P8").Value = 0.4479161786269528
$ws.Range("Q8").Value = 22.94015677745422
$ws.Range("R8").Value = 206.461410997088
$ws.Range("S8").Value = 0.2703659600555272
$ws.Range("T8").Value = 0.2703659600555273

# Row 9: FAPs -> Inflammatory-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Itga8"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.28369633333333
$ws.Range("H9").Value = 30.851089
$ws.Range("I9").Value = 0.6036083824529627
$ws.Range("J9").Value = 0.6036083824529627
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.029162
$ws.Range("N9").Value = 3.087486
$ws.Range("O9").Value = 0.2066490218278579
$ws.Range("P9").Value = 0.2066490218278579
$ws.Range("Q9").Value = 10.583589485806
$ws.Range("R9").Value = 95.25230537225401
$ws.Range("S9").Value = 0.1247350818010003
$ws.Range("T9").Value = 0.1247350818010003

# Row 10: FAPs -> MuSCs
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Itga8"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.28369633333333
$ws.Range("H10").Value = 30.851089
$ws.Range("I10").Value = 0.6036083824529627
$ws.Range("J10").Value = 0.6036083824529627
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8350726666666667
$ws.Range("N10").Value = 2.505218
$ws.Range("O10").Value = 0.1676771487111334
$ws.Range("P10").Value = 0.1676771487111334
$ws.Range("Q10").Value = 8.587633720266888
$ws.Range("R10").Value = 77.28870348240201
$ws.Range("S10").Value = 0.1012113325078521
$ws.Range("T10").Value = 0.1012113325078521

# Row 11: FAPs -> Resolving-Mac
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Vtn"
$ws.Range("C11").Value = "Itga8"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.28369633333333
$ws.Range("H11").Value = 30.851089
$ws.Range("I11").Value = 0.6036083824529627
$ws.Range("J11").Value = 0.6036083824529627
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.5144246666666666
$ws.Range("N11").Value = 1.543274
$ws.Range("O11").Value = 0.1032931201995298
$ws.Range("P11").Value = 0.1032931201995298
$ws.Range("Q11").Value = 5.290187058376222
$ws.Range("R11").Value = 47.61168352538601
$ws.Range("S11").Value = 0.06234859320215763
$ws.Range("T11").Value = 0.06234859320215765

# Row 12: MuSCs -> ECs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Vtn"
$ws.Range("C12").Value = "Itga8"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.355061
$ws.Range("H12").Value = 10.065183
$ws.Range("I12").Value = 0.196927532435664
$ws.Range("J12").Value = 0.196927532435664
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.3708513333333334
$ws.Range("N12").Value = 1.112554
$ws.Range("O12").Value = 0.07446453063452614
$ws.Range("P12").Value = 0.07446453063452615
$ws.Range("Q12").Value = 1.244228845264667
$ws.Range("R12").Value = 11.198059607382
$ws.Range("S12").Value = 0.01466411627183714
$ws.Range("T12").Value = 0.01466411627183715

# Row 13: MuSCs -> FAPs
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Vtn"
$ws.Range("C13").Value = "Itga8"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.355061
$ws.Range("H13").Value = 10.065183
$ws.Range("I13").Value = 0.196927532435664
$ws.Range("J13").Value = 0.196927532435664
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.230730666666667
$ws.Range("N13").Value = 6.692192
$ws.Range("O13").Value = 0.4479161786269527
$ws.Range("P13").Value = 0.4479161786269528
$ws.Range("Q13").Value = 7.484237461237332
$ws.Range("R13").Value = 67.358137151136
$ws.Range("S13").Value = 0.0882070277950179
$ws.Range("T13").Value = 0.08820702779501792

# Row 14: MuSCs -> Inflammatory-Mac
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Vtn"
$ws.Range("C14").Value = "Itga8"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3.355061
$ws.Range("H14").Value = 10.065183
$ws.Range("I14").Value = 0.196927532435664
$ws.Range("J14").Value = 0.196927532435664
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.029162
$ws.Range("N14").Value = 3.087486
$ws.Range("O14").Value = 0.2066490218278579
$ws.Range("P14").Value = 0.2066490218278579
$ws.Range("Q14").Value = 3.452901288882
$ws.Range("R14").Value = 31.076111599938
$ws.Range("S14").Value = 0.04069488194880373
$ws.Range("T14").Value = 0.04069488194880373

# Row 15: MuSCs -> MuSCs
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Vtn"
$ws.Range("C15").Value = "Itga8"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.355061
$ws.Range("H15").Value = 10.065183
$ws.Range("I15").Value = 0.196927532435664
$ws.Range("J15").Value = 0.196927532435664
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.8350726666666667
$ws.Range("N15").Value = 2.505218
$ws.Range("O15").Value = 0.1676771487111334
$ws.Range("P15").Value = 0.1676771487111334
$ws.Range("Q15").Value = 2.801719736099333
$ws.Range("R15").Value = 25.215477624894
$ws.Range("S15").Value = 0.03302024714153138
$ws.Range("T15").Value = 0.03302024714153139

# Row 16: MuSCs -> Resolving-Mac
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Vtn"
$ws.Range("C16").Value = "Itga8"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 3.355061
$ws.Range("H16").Value = 10.065183
$ws.Range("I16").Value = 0.196927532435664
$ws.Range("J16").Value = 0.196927532435664
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.5144246666666666
$ws.Range("N16").Value = 1.543274
$ws.Range("O16").Value = 0.1032931201995298
$ws.Range("P16").Value = 0.1032931201995298
$ws.Range("Q16").Value = 1.725926136571333
$ws.Range("R16").Value = 15.533335229142
$ws.Range("S16").Value = 0.02034125927847385
$ws.Range("T16").Value = 0.02034125927847385

# Row 17: Resolving-Mac -> ECs
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Vtn"
$ws.Range("C17").Value = "Itga8"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.06096333333333333
$ws.Range("H17").Value = 0.18289
$ws.Range("I17").Value = 0.003578283316573439
$ws.Range("J17").Value = 0.003578283316573439
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3708513333333334
$ws.Range("N17").Value = 1.112554
$ws.Range("O17").Value = 0.07446453063452614
$ws.Range("P17").Value = 0.07446453063452615
$ws.Range("Q17").Value = 0.02260833345111112
$ws.Range("R17").Value = 0.20347500106
$ws.Range("S17").Value = 0.0002664551876459966
$ws.Range("T17").Value = 0.0002664551876459967

# Row 18: Resolving-Mac -> FAPs
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Vtn"
$ws.Range("C18").Value = "Itga8"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.06096333333333333
$ws.Range("H18").Value = 0.18289
$ws.Range("I18").Value = 0.003578283316573439
$ws.Range("J18").Value = 0.003578283316573439
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 2.230730666666667
$ws.Range("N18").Value = 6.692192
$ws.Range("O18").Value = 0.4479161786269527
$ws.Range("P18").Value = 0.4479161786269528
$ws.Range("Q18").Value = 0.1359927772088889
$ws.Range("R18").Value = 1.22393499488
$ws.Range("S18").Value = 0.001602770989204153
$ws.Range("T18").Value = 0.001602770989204153

# Row 19: Resolving-Mac -> Inflammatory-Mac
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Vtn"
$ws.Range("C19").Value = "Itga8"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.06096333333333333
$ws.Range("H19").Value = 0.18289
$ws.Range("I19").Value = 0.003578283316573439
$ws.Range("J19").Value = 0.003578283316573439
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 1.029162
$ws.Range("N19").Value = 3.087486
$ws.Range("O19").Value = 0.2066490218278579
$ws.Range("P19").Value = 0.2066490218278579
$ws.Range("Q19").Value = 0.06274114606
$ws.Range("R19").Value = 0.56467031454
$ws.Range("S19").Value = 0.0007394487471928444
$ws.Range("T19").Value = 0.0007394487471928442

# Row 20: Resolving-Mac -> MuSCs
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Vtn"
$ws.Range("C20").Value = "Itga8"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.06096333333333333
$ws.Range("H20").Value = 0.18289
$ws.Range("I20").Value = 0.003578283316573439
$ws.Range("J20").Value = 0.003578283316573439
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.8350726666666667
$ws.Range("N20").Value = 2.505218
$ws.Range("O20").Value = 0.1676771487111334
$ws.Range("P20").Value = 0.1676771487111334
$ws.Range("Q20").Value = 0.05090881333555556
$ws.Range("R20").Value = 0.45817932002
$ws.Range("S20").Value = 0.0005999963438036522
$ws.Range("T20").Value = 0.0005999963438036522

# Row 21: Resolving-Mac -> Resolving-Mac
$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Vtn"
$ws.Range("C21").Value = "Itga8"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.06096333333333333
$ws.Range("H21").Value = 0.18289
$ws.Range("I21").Value = 0.003578283316573439
$ws.Range("J21").Value = 0.003578283316573439
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.5144246666666666
$ws.Range("N21").Value = 1.543274
$ws.Range("O21").Value = 0.1032931201995298
$ws.Range("P21").Value = 0.1032931201995298
$ws.Range("Q21").Value = 0.03136104242888889
$ws.Range("R21").Value = 0.28224938186
$ws.Range("S21").Value = 0.0003696120487267924
$ws.Range("T21").Value = 0.0003696120487267924
